# Append a new time-log entry to the bottom of the log (row 29).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 29

# Column A: date/time entry, formatted like the row above it (style index 3 in the xml == same
# number format / border as A28).
$cellA = $ws.Cells.Item($newRow, 1)
$cellA.Value = "3/13, 5 hrs"

# Column B: work done description, formatted like the row above it (style index 4 == same
# font / border as B28).
$cellB = $ws.Cells.Item($newRow, 2)
$cellB.Value = "Final touches, documenting, recording video, making presentation"

# Copy formatting (borders, font, alignment, number format) from the row above so the new
# row matches the rest of the table, including its wrapped-text row height.
$ws.Range("A28:B28").Copy() | Out-Null
$ws.Range("A29:B29").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# The row above is 41.4pt tall to fit its wrapped text; match it so the new entry looks the
# same as the rest of the log.
$ws.Rows.Item($newRow).RowHeight = 41.4

# Move the active selection to reflect where the user would type the next entry.
$ws.Range("B30").Select() | Out-Null
